# Junction_Flooding_68.xlsx edit
#
# 1) Row 5 values are rewritten at "custom accuracy" (2 decimal places
#    instead of 3).
# 2) Row 6 (the extra 1000th data point that no longer belongs on this
#    sheet) is deleted outright, which also shrinks the used range from
#    A1:AH6 down to A1:AH5.
# 3) A handful of column widths shift by 1 unit (mostly 8 -> 7) because
#    Excel auto-sizes columns to fit the new, shorter numeric strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Re-write row 5 to 2-decimal precision ------------------------------
# (B5, F5 and AE5 already round-trip unchanged at 2/3 decimals, so they are
# left untouched.)
$ws.Range("C5").Value  = 13.21
$ws.Range("D5").Value  = 0.48
$ws.Range("E5").Value  = 38.63
$ws.Range("G5").Value  = 13.65
$ws.Range("H5").Value  = 49.45
$ws.Range("I5").Value  = 21.33
$ws.Range("J5").Value  = 9.65
$ws.Range("K5").Value  = 14.43
$ws.Range("L5").Value  = 15.93
$ws.Range("M5").Value  = 16.54
$ws.Range("N5").Value  = 4.56
$ws.Range("O5").Value  = 13.74
$ws.Range("P5").Value  = 19.84
$ws.Range("Q5").Value  = 11.4
$ws.Range("R5").Value  = 0.15
$ws.Range("S5").Value  = 0.72
$ws.Range("T5").Value  = 204.18
$ws.Range("U5").Value  = 38.43
$ws.Range("V5").Value  = 12.67
$ws.Range("W5").Value  = 26.2
$ws.Range("X5").Value  = 13.64
$ws.Range("Y5").Value  = 2.05
$ws.Range("Z5").Value  = 24.69
$ws.Range("AA5").Value = 11.31
$ws.Range("AB5").Value = 10.13
$ws.Range("AC5").Value = 11.83
$ws.Range("AD5").Value = 16.69
$ws.Range("AF5").Value = 45.01
$ws.Range("AG5").Value = 7.47
$ws.Range("AH5").Value = 15.89

# --- 2) Drop row 6 ----------------------------------------------------------
$ws.Rows.Item(6).Delete()

# --- 3) Shrink the affected columns by one unit -----------------------------
# Excel's ColumnWidth property is offset from the stored OOXML column width
# by 5/6 of a character, so "target width N" is set via ColumnWidth = N - 5/6.
$targetWidths = @{
    3  = 7; 4  = 6; 7  = 7; 9  = 7; 11 = 7; 12 = 7; 13 = 7; 15 = 7; 16 = 7;
    17 = 7; 20 = 8; 22 = 7; 23 = 7; 24 = 7; 26 = 7; 27 = 7; 28 = 7; 29 = 7;
    30 = 7; 34 = 7
}
foreach ($col in $targetWidths.Keys) {
    $ws.Columns.Item($col).ColumnWidth = $targetWidths[$col] - (5/6)
}

Write-Host "edit complete"
